$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update selection to F16 (cosmetic, matches the selected cell in the diff)
$ws.Range("F16").Select()

# E10: 14 -> 8
$ws.Range("E10").Value = 8

# E11: "10:00 AM To 01:00 PM" -> "12:00 AM To 03:00 PM" (new shared string)
$ws.Range("E11").Value = "12:00 AM To 03:00 PM"

# Rows 14-18 and 21: D -> 0, E -> 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0

$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0

$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 0

$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0

$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 0

# Rows 19-20: D stays 3, E -> 3
$ws.Range("E19").Value = 3
$ws.Range("E20").Value = 3

# Row 21: D -> 0, E -> 0
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0
